# Standardized colors and adding wide data
$wb = $excel.ActiveWorkbook

# Rename "Preventative Health" to "Prevention" in the reason column
# (rows 3, 8, 13 of Sheet1 - one entry per location group)
$ws = $wb.Worksheets("Sheet1")
$ws.Range("A3").Value = "Prevention"
$ws.Range("A8").Value = "Prevention"
$ws.Range("A13").Value = "Prevention"

# Remove the now-unused empty Sheet2 and Sheet3
$wb.Worksheets("Sheet2").Delete() | Out-Null
$wb.Worksheets("Sheet3").Delete() | Out-Null
